$wb = $excel.ActiveWorkbook

# ALC row 33: Glazed and Confused / Clear Glass Lens
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 90.1875  # H33: 133.53334 -> 90.1875
$ws.Cells.Item(33, 9).Value = 64.64286  # I33: 107.5 -> 64.64286
$ws.Cells.Item(33, 10).Value = 269  # J33: 498 -> 269
$ws.Cells.Item(33, 11).Value = 64.64286  # K33: 107.5 -> 64.64286
$ws.Cells.Item(33, 12).Value = 269  # L33: 498 -> 269
$ws.Cells.Item(33, 13).Value = 164.35714  # M33: 121.5 -> 164.35714
$ws.Cells.Item(33, 14).Value = -727  # N33: -956 -> -727

# ALC row 55: A Real Smooth Move / Lanolin
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(55, 8).Value = 116.666664  # H55: 137 -> 116.666664
$ws.Cells.Item(55, 9).Value = 150  # I55: 137 -> 150
$ws.Cells.Item(55, 10).Value = 100  # J55: 0 -> 100
$ws.Cells.Item(55, 11).Value = 150  # K55: 137 -> 150
$ws.Cells.Item(55, 12).Value = 100  # L55: 0 -> 100
$ws.Cells.Item(55, 13).Value = 64  # M55: 77 -> 64
$ws.Cells.Item(55, 14).Value = -528  # N55: (new cell) -> -528

# ALC row 94: Magic Beans / Growth Formula Eta
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(94, 8).Value = 998  # H94: 909.5 -> 998
$ws.Cells.Item(94, 9).Value = 998  # I94: 909.5 -> 998
$ws.Cells.Item(94, 11).Value = 998  # K94: 909.5 -> 998
$ws.Cells.Item(94, 13).Value = -547  # M94: -458.5 -> -547

# ALC row 96: Scroll Down / Grade 1 Reisui of Intelligence
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(96, 8).Value = 884.375  # H96: 901.875 -> 884.375
$ws.Cells.Item(96, 9).Value = 725.1429000000001  # I96: 745.1429000000001 -> 725.1429000000001
$ws.Cells.Item(96, 11).Value = 2175.4287  # K96: 2235.4287 -> 2175.4287
$ws.Cells.Item(96, 13).Value = -802.4287000000004  # M96: -862.4287000000004 -> -802.4287000000004

# ALC row 98: The Dotted Line / Enchanted Durium Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 716.1539  # H98: 734.1539 -> 716.1539
$ws.Cells.Item(98, 9).Value = 351  # I98: 384.42856 -> 351
$ws.Cells.Item(98, 11).Value = 351  # K98: 384.42856 -> 351
$ws.Cells.Item(98, 13).Value = 1147  # M98: 1113.57144 -> 1147

# ALC row 111: An Eye for Healing / Grade 1 Dexterity Alkahest
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(111, 8).Value = 4313.9  # H111: 5084.222 -> 4313.9
$ws.Cells.Item(111, 9).Value = 4283.8887  # I111: 5084.222 -> 4283.8887
$ws.Cells.Item(111, 10).Value = 4584  # J111: 0 -> 4584
$ws.Cells.Item(111, 11).Value = 12851.6661  # K111: 15252.666 -> 12851.6661
$ws.Cells.Item(111, 12).Value = 13752  # L111: 0 -> 13752
$ws.Cells.Item(111, 13).Value = -9784.666100000002  # M111: -12185.666 -> -9784.666100000002
$ws.Cells.Item(111, 14).Value = -19886  # N111: (new cell) -> -19886

# ALC row 122: Wishful Inking / Enchanted High Durium Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 716.1539  # H122: 734.1539 -> 716.1539
$ws.Cells.Item(122, 9).Value = 351  # I122: 384.42856 -> 351
$ws.Cells.Item(122, 11).Value = 1053  # K122: 1153.28568 -> 1053
$ws.Cells.Item(122, 13).Value = 1397  # M122: 1296.71432 -> 1397

# ALC row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 2175.6843  # H137: 2209.2163 -> 2175.6843
$ws.Cells.Item(137, 9).Value = 1340.3572  # I137: 1355.3704 -> 1340.3572
$ws.Cells.Item(137, 11).Value = 4021.0716  # K137: 4066.1112 -> 4021.0716
$ws.Cells.Item(137, 13).Value = -1471.0716  # M137: -1516.1112 -> -1471.0716

# ALC row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2384.328  # H138: 2408.0327 -> 2384.328
$ws.Cells.Item(138, 10).Value = 2615.1956  # J138: 2646.6304 -> 2615.1956
$ws.Cells.Item(138, 12).Value = 7845.5868  # L138: 7939.8912 -> 7845.5868
$ws.Cells.Item(138, 14).Value = -18125.5868  # N138: -18219.8912 -> -18125.5868

# ARM row 5: The Alloyed Truth / Bronze Rivets
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 192.5  # H5: 134.5 -> 192.5
$ws.Cells.Item(5, 9).Value = 90  # I5: 100 -> 90
$ws.Cells.Item(5, 10).Value = 500  # J5: 140.25 -> 500
$ws.Cells.Item(5, 11).Value = 90  # K5: 100 -> 90
$ws.Cells.Item(5, 12).Value = 500  # L5: 140.25 -> 500
$ws.Cells.Item(5, 13).Value = 22  # M5: 12 -> 22
$ws.Cells.Item(5, 14).Value = -724  # N5: -364.25 -> -724

# ARM row 32: Ingot We Trust / Steel Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4277.885  # H32: 4277.923 -> 4277.885
$ws.Cells.Item(32, 9).Value = 4289.04  # I32: 4289.08 -> 4289.04
$ws.Cells.Item(32, 11).Value = 4289.04  # K32: 4289.08 -> 4289.04
$ws.Cells.Item(32, 13).Value = -4002.04  # M32: -4002.08 -> -4002.04

# ARM row 74: As the Bolt Flies / Titanium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 2366.3333  # H74: 1670.3334 -> 2366.3333
$ws.Cells.Item(74, 9).Value = 2366.3333  # I74: 1670.3334 -> 2366.3333
$ws.Cells.Item(74, 11).Value = 2366.3333  # K74: 1670.3334 -> 2366.3333
$ws.Cells.Item(74, 13).Value = -1492.3333  # M74: -796.3334 -> -1492.3333

# ARM row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 2366.3333  # H77: 1670.3334 -> 2366.3333
$ws.Cells.Item(77, 9).Value = 2366.3333  # I77: 1670.3334 -> 2366.3333
$ws.Cells.Item(77, 11).Value = 11831.6665  # K77: 8351.666999999999 -> 11831.6665
$ws.Cells.Item(77, 13).Value = -7463.666499999999  # M77: -3983.666999999999 -> -7463.666499999999

# ARM row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 1666  # H102: 1921.7778 -> 1666
$ws.Cells.Item(102, 9).Value = 832  # I102: 1459.4 -> 832
$ws.Cells.Item(102, 10).Value = 2500  # J102: 2499.75 -> 2500
$ws.Cells.Item(102, 11).Value = 832  # K102: 1459.4 -> 832
$ws.Cells.Item(102, 12).Value = 2500  # L102: 2499.75 -> 2500
$ws.Cells.Item(102, 13).Value = 790  # M102: 162.5999999999999 -> 790
$ws.Cells.Item(102, 14).Value = -5744  # N102: -5743.75 -> -5744

# BSM row 4: Mending Fences / Bronze Rivets
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 192.5  # H4: 134.5 -> 192.5
$ws.Cells.Item(4, 9).Value = 90  # I4: 100 -> 90
$ws.Cells.Item(4, 10).Value = 500  # J4: 140.25 -> 500
$ws.Cells.Item(4, 11).Value = 90  # K4: 100 -> 90
$ws.Cells.Item(4, 12).Value = 500  # L4: 140.25 -> 500
$ws.Cells.Item(4, 13).Value = 25  # M4: 15 -> 25
$ws.Cells.Item(4, 14).Value = -730  # N4: -370.25 -> -730

# BSM row 44: You Spin Me Round / Mythril Broadsword
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(44, 8).Value = 0  # H44: 5045 -> 0
$ws.Cells.Item(44, 9).Value = 0  # I44: 5045 -> 0
$ws.Cells.Item(44, 11).Value = 0  # K44: 5045 -> 0
$ws.Cells.Item(44, 13).ClearContents()  # M44: -4548 -> (removed)

# BSM row 105: Ingot to Wing It / Molybdenum Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 4492.25  # H105: 5058.5 -> 4492.25
$ws.Cells.Item(105, 9).Value = 4492.25  # I105: 5058.5 -> 4492.25
$ws.Cells.Item(105, 11).Value = 4492.25  # K105: 5058.5 -> 4492.25
$ws.Cells.Item(105, 13).Value = -2745.25  # M105: -3311.5 -> -2745.25

# BSM row 133: Paring Is Caring / Mountain Chromite Hatchet
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(133, 8).Value = 74998.336  # H133: 74998.75 -> 74998.336
$ws.Cells.Item(133, 10).Value = 74998.336  # J133: 74998.75 -> 74998.336
$ws.Cells.Item(133, 12).Value = 74998.336  # L133: 74998.75 -> 74998.336
$ws.Cells.Item(133, 14).Value = -85118.336  # N133: -85118.75 -> -85118.336

# CRP row 31: Wall Not Found / Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2130.125  # H31: 2220.5 -> 2130.125
$ws.Cells.Item(31, 9).Value = 2228.3076  # I31: 2361.182 -> 2228.3076
$ws.Cells.Item(31, 11).Value = 2228.3076  # K31: 2361.182 -> 2228.3076
$ws.Cells.Item(31, 13).Value = -1933.3076  # M31: -2066.182 -> -1933.3076

# CRP row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2130.125  # H34: 2220.5 -> 2130.125
$ws.Cells.Item(34, 9).Value = 2228.3076  # I34: 2361.182 -> 2228.3076
$ws.Cells.Item(34, 11).Value = 2228.3076  # K34: 2361.182 -> 2228.3076
$ws.Cells.Item(34, 13).Value = -2026.3076  # M34: -2159.182 -> -2026.3076

# CRP row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2060.9092  # H58: 2324.8572 -> 2060.9092
$ws.Cells.Item(58, 9).Value = 1677.8572  # I58: 1783 -> 1677.8572
$ws.Cells.Item(58, 11).Value = 1677.8572  # K58: 1783 -> 1677.8572
$ws.Cells.Item(58, 13).Value = -1474.8572  # M58: -1580 -> -1474.8572

# CRP row 86: Birch, Please / Birch Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 7787.643  # H86: 8502.833000000001 -> 7787.643
$ws.Cells.Item(86, 9).Value = 8488.777  # I86: 11601.2 -> 8488.777
$ws.Cells.Item(86, 10).Value = 6525.6  # J86: 6289.7144 -> 6525.6
$ws.Cells.Item(86, 11).Value = 8488.777  # K86: 11601.2 -> 8488.777
$ws.Cells.Item(86, 12).Value = 6525.6  # L86: 6289.7144 -> 6525.6
$ws.Cells.Item(86, 13).Value = -7365.777  # M86: -10478.2 -> -7365.777
$ws.Cells.Item(86, 14).Value = -8771.6  # N86: -8535.714400000001 -> -8771.6

# CRP row 89: Built This City on Blocks and Soul (L) / Birch Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 7787.643  # H89: 8502.833000000001 -> 7787.643
$ws.Cells.Item(89, 9).Value = 8488.777  # I89: 11601.2 -> 8488.777
$ws.Cells.Item(89, 10).Value = 6525.6  # J89: 6289.7144 -> 6525.6
$ws.Cells.Item(89, 11).Value = 42443.885  # K89: 58006 -> 42443.885
$ws.Cells.Item(89, 12).Value = 32628  # L89: 31448.572 -> 32628
$ws.Cells.Item(89, 13).Value = -36827.885  # M89: -52390 -> -36827.885
$ws.Cells.Item(89, 14).Value = -43860  # N89: -42680.572 -> -43860

# CRP row 94: Beech, Please / Beech Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value = 833.3333  # H94: 674.75 -> 833.3333
$ws.Cells.Item(94, 10).Value = 250  # J94: 233 -> 250
$ws.Cells.Item(94, 12).Value = 250  # L94: 233 -> 250
$ws.Cells.Item(94, 14).Value = -1152  # N94: -1135 -> -1152

# CRP row 107: Built to Last / White Oak Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 1454.9166  # H107: 1359.2307 -> 1454.9166
$ws.Cells.Item(107, 9).Value = 617.2857  # I107: 566.5 -> 617.2857
$ws.Cells.Item(107, 11).Value = 617.2857  # K107: 566.5 -> 617.2857
$ws.Cells.Item(107, 13).Value = 1302.7143  # M107: 1353.5 -> 1302.7143

# CRP row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 4028.3333  # H122: 4070.875 -> 4028.3333
$ws.Cells.Item(122, 9).Value = 3759  # I122: 3773.2 -> 3759
$ws.Cells.Item(122, 11).Value = 11277  # K122: 11319.6 -> 11277
$ws.Cells.Item(122, 13).Value = -8827  # M122: -8869.599999999999 -> -8827

# CRP row 136: Turali Quality / Dark Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 2060.9092  # H136: 2324.8572 -> 2060.9092
$ws.Cells.Item(136, 9).Value = 1677.8572  # I136: 1783 -> 1677.8572
$ws.Cells.Item(136, 11).Value = 5033.571599999999  # K136: 5349 -> 5033.571599999999
$ws.Cells.Item(136, 13).Value = -2483.571599999999  # M136: -2799 -> -2483.571599999999

# CUL row 137: Creative Chocolate / Gateau au Chocolat
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 8).Value = 2485.5715  # H137: 2822 -> 2485.5715
$ws.Cells.Item(137, 10).Value = 4000  # J137: 3999.6667 -> 4000
$ws.Cells.Item(137, 12).Value = 12000  # L137: 11999.0001 -> 12000
$ws.Cells.Item(137, 14).Value = -22200  # N137: -22199.0001 -> -22200

# GSM row 70: Sky Is the Limit / Mythrite Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 9).Value = 0  # I70: 6000 -> 0
$ws.Cells.Item(70, 11).Value = 0  # K70: 6000 -> 0
$ws.Cells.Item(70, 13).ClearContents()  # M70: -5730 -> (removed)

# GSM row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 9).Value = 0  # I73: 6000 -> 0
$ws.Cells.Item(73, 11).Value = 0  # K73: 6000 -> 0
$ws.Cells.Item(73, 13).ClearContents()  # M73: -5064 -> (removed)

# GSM row 80: Needs More Prayerbell / Hardsilver Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3609.7  # H80: 3554.182 -> 3609.7
$ws.Cells.Item(80, 9).Value = 2898.7144  # I80: 2911.25 -> 2898.7144
$ws.Cells.Item(80, 11).Value = 2898.7144  # K80: 2911.25 -> 2898.7144
$ws.Cells.Item(80, 13).Value = -1900.7144  # M80: -1913.25 -> -1900.7144

# GSM row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 3609.7  # H83: 3554.182 -> 3609.7
$ws.Cells.Item(83, 9).Value = 2898.7144  # I83: 2911.25 -> 2898.7144
$ws.Cells.Item(83, 11).Value = 14493.572  # K83: 14556.25 -> 14493.572
$ws.Cells.Item(83, 13).Value = -9501.572  # M83: -9564.25 -> -9501.572

# GSM row 102: Put the Metal to the Peddle / Durium Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 6482.8  # H102: 6978.5 -> 6482.8
$ws.Cells.Item(102, 9).Value = 4800  # I102: 4950 -> 4800
$ws.Cells.Item(102, 11).Value = 4800  # K102: 4950 -> 4800
$ws.Cells.Item(102, 13).Value = -3178  # M102: -3328 -> -3178

# GSM row 113: Copious Crystal Cannons / Manasilver Nugget
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1986.8572  # H113: 2151.8333 -> 1986.8572
$ws.Cells.Item(113, 9).Value = 1331.6666  # I113: 1499 -> 1331.6666
$ws.Cells.Item(113, 11).Value = 1331.6666  # K113: 1499 -> 1331.6666
$ws.Cells.Item(113, 13).Value = 838.3334  # M113: 671 -> 838.3334

# WVR row 119: A Job Well Done / Dwarven Cotton Gaskins of Fending
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119, 8).Value = 0  # H119: 70000 -> 0
$ws.Cells.Item(119, 10).Value = 0  # J119: 70000 -> 0
$ws.Cells.Item(119, 12).ClearContents()  # L119: 70000 -> (removed)
$ws.Cells.Item(119, 14).Value = 0  # N119: -79676 -> 0
